# Weekly update: insert the newest week's pricing rows (Primera/Segunda)
# at the top of the "Apio" dataset, pushing the existing rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 48:49 (rows 48.. shift down to 50..)
$ws.Rows("48:49").Insert()

# Row 48 - new "Primera" quality entry for the latest week
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44525
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 100112017
$ws.Range("G48").Value = "Apio"
$ws.Range("H48").Value = "Americana (o)"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 61
$ws.Range("K48").Value = 7000
$ws.Range("L48").Value = 8000
$ws.Range("M48").Value = 7508
$ws.Range("N48").Value = "$/docena de matas"
$ws.Range("O48").Value = "Región Metropolitana"
$ws.Range("P48").Value = 1251
$ws.Range("Q48").Value = 6
$ws.Range("R48").Value = "Hortaliza"

# Row 49 - new "Segunda" quality entry for the latest week
$ws.Range("A49").Value = 9
$ws.Range("B49").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 44525
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = 100112017
$ws.Range("G49").Value = "Apio"
$ws.Range("H49").Value = "Americana (o)"
$ws.Range("I49").Value = "Segunda"
$ws.Range("J49").Value = 34
$ws.Range("K49").Value = 6000
$ws.Range("L49").Value = 6000
$ws.Range("M49").Value = 6000
$ws.Range("N49").Value = "$/docena de matas"
$ws.Range("O49").Value = "Región Metropolitana"
$ws.Range("P49").Value = 1000
$ws.Range("Q49").Value = 6
$ws.Range("R49").Value = "Hortaliza"
